# Insert a new row for "CJB" (Coimbatore, India) right above the "IAD"
# (Ashburn) row, i.e. at row 275, pushing every subsequent data row down
# by one. The sheet grows from A1:H333 to A1:H334.
#
# We avoid Rows.Insert() because (in this runtime) it auto-generates a
# fresh cellXf for the newly inserted row's formatting, which would leave
# an extra/unused style in styles.xml that the real edit never produced.
# Instead we manually extend the used range by copying the last row's
# formatting down, then shift every row's values down one at a time
# (bottom-up, values-only) and overwrite row 275 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the grid by one row, copying row 333's formatting onto the new
# row 334 so the shifted-down content (originally row 333, "YHZ") keeps
# its existing look (bold/bordered colo-code cell, etc.).
$ws.Range("A333:H333").Copy()
$ws.Range("A334:H334").PasteSpecial(-4122)

# Shift rows 275..333 down into 276..334, working from the bottom up so
# nothing gets clobbered before it's copied.
for ($r = 333; $r -ge 275; $r--) {
    $src = $ws.Range("A" + $r + ":H" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":H" + ($r + 1))
    $src.Copy()
    $dst.PasteSpecial(-4104)
}

# Write the new Coimbatore, India record into the now-vacated row 275.
$ws.Range("A275").Value = "CJB"
$ws.Range("B275").Value = "Coimbatore, India"
$ws.Range("C275").Value = "Asia Pacific"
$ws.Range("D275").Value = "Coimbatore"
$ws.Range("E275").Value = "India"
$ws.Range("F275").Value = "IN"
$ws.Range("G275").Value = 11.008756
$ws.Range("H275").Value = 77.011987
